$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: counts (episode count=40 unchanged, arrival count, collision count, local minima count)
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 28

# Row 5: rates (arrival rate, collision rate, minima rate)
$ws.Range("B5").Value = 0.3
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.7
